# Applies updated crypto price/volume figures to Sheet1 (columns D and E),
# matching the "Updated cryptos list" GitHub Actions commit.
#
# The Price/Volume columns hold plain text in the source data (e.g. prices
# such as "74.927.63" use dots as thousands separators, and percentages are
# padded with spaces). For values that look like ordinary numbers (e.g.
# "1.00", "0.194"), a leading apostrophe is used so Excel keeps them as text
# instead of coercing them into a Double and silently dropping significant
# trailing zeros (e.g. "1.00" -> 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '75.004.08'
$ws.Range("E2").Value = '  +6.97%  '
$ws.Range("D3").Value = '2.674.37'
$ws.Range("E3").Value = '  +9.08%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'186.89"
$ws.Range("E5").Value = '  +11.88%  '
$ws.Range("D6").Value = "'586.34"
$ws.Range("E6").Value = '  +3.08%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  +3.88%  '
$ws.Range("D9").Value = "'0.194"
$ws.Range("E9").Value = '  +10.10%  '
$ws.Range("D10").Value = '2.672.77'
$ws.Range("E10").Value = '  +9.11%  '
$ws.Range("D11").Value = "'0.164"
$ws.Range("E11").Value = '  +1.59%  '
$ws.Range("E12").Value = '  +5.95%  '
$ws.Range("E13").Value = '  +0.50%  '
$ws.Range("D14").Value = '3.162.64'
$ws.Range("E14").Value = '  +8.92%  '
$ws.Range("D15").Value = '74.840.14'
$ws.Range("E15").Value = '  +6.90%  '
$ws.Range("E16").Value = '  +2.61%  '
$ws.Range("D17").Value = "'26.50"
$ws.Range("E17").Value = '  +9.39%  '
$ws.Range("D18").Value = '2.682.23'
$ws.Range("E18").Value = '  +9.08%  '
$ws.Range("D19").Value = "'9.14"
$ws.Range("E19").Value = '  +27.81%  '
$ws.Range("E20").Value = '  +8.85%  '
$ws.Range("D21").Value = "'371.37"
$ws.Range("E21").Value = '  +8.69%  '
$ws.Range("E22").Value = '  +11.94%  '
$ws.Range("E23").Value = '  +4.48%  '
$ws.Range("D24").Value = "'6.26"
$ws.Range("E24").Value = '  +3.64%  '
$ws.Range("E25").Value = '  +0.39%  '
$ws.Range("D26").Value = "'69.69"
$ws.Range("E26").Value = '  +4.82%  '
$ws.Range("D27").Value = "'4.14"
$ws.Range("E27").Value = '  +7.69%  '
$ws.Range("D28").Value = "'9.34"
$ws.Range("E28").Value = '  +9.28%  '
$ws.Range("D29").Value = '2.795.14'
$ws.Range("E29").Value = '  +8.31%  '
$ws.Range("E30").Value = '  +2.80%  '
$ws.Range("E31").Value = '  +10.33%  '
$ws.Range("E32").Value = '  +13.38%  '
$ws.Range("D33").Value = "'521.92"
$ws.Range("E33").Value = '  +13.16%  '
$ws.Range("D34").Value = "'7.65"
$ws.Range("E34").Value = '  +3.48%  '
$ws.Range("E35").Value = '  +7.74%  '
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("D37").Value = "'163.53"
$ws.Range("E37").Value = '  +1.85%  '
$ws.Range("E38").Value = '  +5.24%  '
$ws.Range("D39").Value = "'19.19"
$ws.Range("E39").Value = '  +5.38%  '
$ws.Range("D40").Value = "'19.33"
$ws.Range("E40").Value = '  +1.13%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("E42").Value = '  +12.85%  '
$ws.Range("D43").Value = "'169.66"
$ws.Range("E43").Value = '  +26.11%  '
$ws.Range("D44").Value = "'0.330"
$ws.Range("E44").Value = '  +8.55%  '
$ws.Range("E45").Value = '  +9.19%  '
$ws.Range("E46").Value = '  +8.47%  '
$ws.Range("E47").Value = '  +10.84%  '
$ws.Range("D48").Value = "'39.03"
$ws.Range("E48").Value = '  +2.44%  '
$ws.Range("D49").Value = "'0.0841"
$ws.Range("E49").Value = '  +15.65%  '
$ws.Range("D50").Value = "'3.64"
$ws.Range("E50").Value = '  +6.96%  '
$ws.Range("D51").Value = "'0.529"
$ws.Range("E51").Value = '  +7.69%  '
